$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 4233.283585129627
    3  = 4161.36506566689
    4  = 4161.36506566689
    5  = 4161.36506566689
    6  = 3979.461954550371
    7  = 3874.302871801664
    8  = 3874.302871801664
    9  = 3874.302871801664
    10 = 3874.302871801664
    11 = 3868.628057416159
    12 = 3857.224278352653
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 3).Value = $values[$row]
}
